# Fill in the "Method Inputs" / "Expected Result" test-plan data and the
# student name, matching the commit "Tested functionality of the methods
# to see if they perform as expected."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Test cases for __init__ (rows 7-10): Method Inputs (F) filled in,
#     Condition inputs (E) all "None", Expected Result (G) filled in.
#     Entered column-by-column (E then F then G) to mirror the authoring
#     order reflected in the shared-string table.
$ws.Range("E7").Value  = "None"
$ws.Range("E8").Value  = "None"
$ws.Range("E9").Value  = "None"
$ws.Range("E10").Value = "None"

$ws.Range("F7").Value  = '("DUNE", "Frank Herbert", Genre.FICTION)'
$ws.Range("F8").Value  = '("", "Frank Herbert", Genre.FICTION)'
$ws.Range("F9").Value  = '("DUNE", "", Genre.FICTION)'
$ws.Range("F10").Value = '("DUNE", "Frank Herbert", "INVALID)'

$ws.Range("G7").Value  = "Pass/No errors"
$ws.Range("G8").Value  = "Raise ValueError"
$ws.Range("G9").Value  = "Raise ValueError"
$ws.Range("G10").Value = "Raise ValueError"

# --- Test cases for the title/author/genre getters (rows 11-13).
$ws.Range("E11").Value = 'LibraryItem("DUNE", "Frank Herbert", Genre.FICTION)'
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "Returns title"

$ws.Range("E12").Value = 'LibraryItem("DUNE", "Frank Herbert", Genre.FICTION)'
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "Returns author"

$ws.Range("E13").Value = 'LibraryItem("DUNE", "Frank Herbert", Genre.FICTION)'
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "Returns genre"

# --- Student name (merged C3:D3), entered last.
$ws.Range("C3").Value = "Philip Pacla-on"

# --- Page orientation switched to portrait for printing the completed plan.
$ws.PageSetup.Orientation = 1

# --- Leave the active selection on the student-name cell, like the author did.
[void]$ws.Range("C3:D3").Select()
